# GradedExercise - grade entry pass.
# Fills in the "Value" column (G) for every graded criterion, re-applies the
# "Good" cell style to the two rubric rows whose highlight changed, and moves
# the live selection to reflect where the grader ended up (I21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grade values entered into column G --------------------------------
$ws.Range("G4").Value  = 5
$ws.Range("G5").Value  = 15

$ws.Range("G8").Value  = 2
$ws.Range("G9").Value  = 2
$ws.Range("G10").Value = 3
$ws.Range("G12").Value = 5

$ws.Range("G16").Value = 5
$ws.Range("G17").Value = 5
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 2

$ws.Range("G25").Value = 2

# --- Rubric rows re-highlighted from "Neutral" to "Good" ----------------
$ws.Range("D4:E4").Style   = "Good"
$ws.Range("D10:E10").Style = "Good"

# --- Recalculate so every dependent formula carries a fresh cached value -
$excel.Calculate()

# --- Restore the grader's on-screen selection ----------------------------
$ws.Range("I21").Select()
